$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.603.66"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").Value = "1.697.31"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'317.11"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.3946"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").Value = "'0.4017"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "'1.527"
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("D10").Value = "'1.000"
$ws.Range("D11").Value = "'52.68"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'0.08765"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'7.218"
$ws.Range("E13").Value = "  +6.58%  "
$ws.Range("D14").Value = "'23.29"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").Value = "'8.153"
$ws.Range("E15").Value = "  +11.91%  "
$ws.Range("D16").Value = "'0.00001316"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "1.692.85"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "'99.81"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "'0.07071"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'19.70"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("D21").Value = "'6.935"
$ws.Range("E21").Value = "  +4.59%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'14.18"
$ws.Range("D24").Value = "24.608.63"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "'3.137"
$ws.Range("E25").Value = "  +10.88%  "
$ws.Range("D26").Value = "'2.334"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'22.78"
$ws.Range("E27").Value = "  +4.83%  "
$ws.Range("D28").Value = "'162.86"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "'136.61"
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("D30").Value = "'5.193"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").Value = "'7.486"
$ws.Range("E31").Value = "  +10.20%  "
$ws.Range("D32").Value = "1.877.96"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").Value = "'1.089"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").Value = "'0.08588"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").Value = "'7.152"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("D36").Value = "'11.52"
$ws.Range("E36").Value = "  +10.08%  "
$ws.Range("D37").Value = "'0.2741"
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("D38").Value = "'1.931"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'14.51"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "'0.09133"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("D41").Value = "'0.02742"
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("D42").Value = "'1.480"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "'0.7672"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").Value = "'0.7181"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "'15.57"
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("D46").Value = "'2.552"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("D47").Value = "'4.219"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'140.93"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").Value = "'1.327"
$ws.Range("E50").Value = "  +9.13%  "
$ws.Range("D51").Value = "'0.07991"
$ws.Range("E51").Value = "  +2.61%  "
